$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("normal")

# ---------------------------------------------------------------------------
# 1) Workbook window position (cosmetic bookViews/workbookView attribute)
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 6460
$win.Top = 3200

# ---------------------------------------------------------------------------
# 2) New header columns Z / AA / AB on sheet "normal"
#    (written in the same order the new shared-strings entries were
#    appended in the source commit: Z, then AB, then AA)
# ---------------------------------------------------------------------------
$ws.Range("Z1").Value = "INT_upgradeCoinNeed"
$ws.Range("AB1").Value = "INT_upgradeTimeSecondsNeed"
$ws.Range("AA1").Value = "INT_upgradeTechPointNeed"

# give the new header cells the same look as the rest of row 1 (fill + thin
# border all around, centered) by copying format from the existing header
$hdr = $ws.Range("Y1")
$hdr.Copy()
$ws.Range("Z1:AB1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Data for the 24 tiers (3 rows per unit: +50 / +100 / +200 coin-need),
#    repeating every 3 rows from row 2 through row 25.
# ---------------------------------------------------------------------------
$vals = @(50, 100, 200)

for ($r = 2; $r -le 25; $r++) {
    $grp = ($r - 2) % 3
    $v = $vals[$grp]

    $ws.Cells.Item($r, 26).Value = $v   # Z
    $ws.Cells.Item($r, 27).Value = $v   # AA
    $ws.Cells.Item($r, 28).Value = $v   # AB

    $zCell = $ws.Cells.Item($r, 26)
    $aaCell = $ws.Cells.Item($r, 27)
    $abCell = $ws.Cells.Item($r, 28)

    # AA is always the "full thin border" style, centered
    $aaCell.Borders.Item(7).LineStyle = 1
    $aaCell.Borders.Item(10).LineStyle = 1
    $aaCell.Borders.Item(8).LineStyle = 1
    $aaCell.Borders.Item(9).LineStyle = 1
    $aaCell.HorizontalAlignment = -4108
    $aaCell.VerticalAlignment = -4108

    if ($grp -eq 0) {
        # first row of the trio: full thin border on all sides, centered
        foreach ($c in @($zCell, $abCell)) {
            $c.Borders.Item(7).LineStyle = 1
            $c.Borders.Item(10).LineStyle = 1
            $c.Borders.Item(8).LineStyle = 1
            $c.Borders.Item(9).LineStyle = 1
            $c.HorizontalAlignment = -4108
            $c.VerticalAlignment = -4108
        }
    } else {
        # continuation rows: thin border on left/right/bottom only (no top),
        # centered - visually attaches to the row above
        foreach ($c in @($zCell, $abCell)) {
            $c.Borders.Item(7).LineStyle = 1
            $c.Borders.Item(10).LineStyle = 1
            $c.Borders.Item(9).LineStyle = 1
            $c.Borders.Item(8).LineStyle = -4142
            $c.HorizontalAlignment = -4108
            $c.VerticalAlignment = -4108
        }
    }
}

# ---------------------------------------------------------------------------
# 4) Sheet view: scroll back to the top-left corner and move the selection
# ---------------------------------------------------------------------------
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("D14").Select()
